# Thrust_Against_Mass_Calculations.xlsx update
# - add two new worksheets: "Exp_data" and "PID_Values" (after Sheet1)
# - populate them with experimental / PID-tuning data
# - make "Exp_data" the active sheet, at 95% zoom (same zoom applied to all sheets)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- create the two new sheets, in order, right after Sheet1 --------------
$wsExp = $wb.Worksheets.Add($null, $ws1)
$wsExp.Name = "Exp_data"

$wsPid = $wb.Worksheets.Add($null, $wsExp)
$wsPid.Name = "PID_Values"

# --- Exp_data contents ------------------------------------------------------
$wsExp.Range("A1").Value = "Length of pod"
$wsExp.Range("B1").Value = "Perching time (minutes)"
$wsExp.Range("C1").Value = "Hover time (seconds)"
$wsExp.Range("D1").Value = "Oscillation amplitude (0 deg, aluminium plate)"
$wsExp.Range("E1").Value = "Max angle (flat aluminium plate)"

$wsExp.Range("A2").Value = 1
$wsExp.Range("B2").Value = 50

$wsExp.Range("A3").Value = 2
$wsExp.Range("B3").Value = 63

$wsExp.Range("A4").Value = 3
$wsExp.Range("B4").Value = 63

$wsExp.Range("A5").Value = 4
$wsExp.Range("B5").Value = 59

$wsExp.Range("A6").Value = 6

$wsExp.Range("A7").Value = 8
$wsExp.Range("B7").Value = 45

$wsExp.Range("A8").Value = 12
$wsExp.Range("B8").Value = 40

$wsExp.Range("A9").Value = 16
$wsExp.Range("B9").Value = "very difficult to keep stable"

$wsExp.Range("A10").Value = 20

$wsExp.Range("A11").Value = 24

$wsExp.Range("B12").Value = "Should perform these tests at least two times\ each"

# --- PID_Values contents -----------------------------------------------------
$wsPid.Range("A1").Value = "Pod length (cm)"
$wsPid.Range("B1").Value = "A_pkp"
$wsPid.Range("C1").Value = "A_pki"
$wsPid.Range("D1").Value = "A_pkd"
$wsPid.Range("E1").Value = "A_rkp"
$wsPid.Range("F1").Value = "A_rki"
$wsPid.Range("G1").Value = "a_rkd"
$wsPid.Range("H1").Value = "a_ykp"
$wsPid.Range("I1").Value = "a_yki"
$wsPid.Range("J1").Value = "a_ykd"
$wsPid.Range("K1").Value = "r_pkp"
$wsPid.Range("L1").Value = "r_pki"
$wsPid.Range("M1").Value = "r_pkd"
$wsPid.Range("N1").Value = "r_rkp"
$wsPid.Range("O1").Value = "r_rki"
$wsPid.Range("P1").Value = "r_rkd"
$wsPid.Range("Q1").Value = "r_ykp"
$wsPid.Range("R1").Value = "r_yki"
$wsPid.Range("S1").Value = "r_ykd"

$wsPid.Range("A2").Value = 1
$wsPid.Range("B2").Value = 2.3
$wsPid.Range("C2").Value = 0
$wsPid.Range("D2").Value = 0.05
$wsPid.Range("E2").Value = 3.1
$wsPid.Range("F2").Value = 0
$wsPid.Range("G2").Value = 0.05
$wsPid.Range("H2").Value = 9
$wsPid.Range("I2").Value = 0.1
$wsPid.Range("J2").Value = 0
$wsPid.Range("K2").Value = 1500
$wsPid.Range("L2").Value = 30
$wsPid.Range("M2").Value = 0.01
$wsPid.Range("N2").Value = 1050
$wsPid.Range("O2").Value = 35
$wsPid.Range("P2").Value = 0.03
$wsPid.Range("Q2").Value = 1000
$wsPid.Range("R2").Value = 30
$wsPid.Range("S2").Value = 0

$wsPid.Range("A3").Value = 2
$wsPid.Range("B3").Value = 2.5
$wsPid.Range("C3").Value = 0
$wsPid.Range("D3").Value = 0.05
$wsPid.Range("E3").Value = 2.5
$wsPid.Range("F3").Value = 0
$wsPid.Range("G3").Value = 0.05
$wsPid.Range("H3").Value = 9
$wsPid.Range("I3").Value = 0.1
$wsPid.Range("J3").Value = 0
$wsPid.Range("K3").Value = 2500
$wsPid.Range("L3").Value = 30
$wsPid.Range("M3").Value = 0.01
$wsPid.Range("N3").Value = 2500
$wsPid.Range("O3").Value = 30
$wsPid.Range("P3").Value = 0.01
$wsPid.Range("Q3").Value = 1000
$wsPid.Range("R3").Value = 30
$wsPid.Range("S3").Value = 0

$wsPid.Range("A4").Value = 3
$wsPid.Range("B4").Value = 2.7
$wsPid.Range("C4").Value = 0
$wsPid.Range("D4").Value = 0.05
$wsPid.Range("E4").Value = 2.7
$wsPid.Range("F4").Value = 0
$wsPid.Range("G4").Value = 0.05
$wsPid.Range("H4").Value = 9
$wsPid.Range("I4").Value = 0.1
$wsPid.Range("J4").Value = 0
$wsPid.Range("K4").Value = 2600
$wsPid.Range("L4").Value = 30
$wsPid.Range("M4").Value = 0.01
$wsPid.Range("N4").Value = 2600
$wsPid.Range("O4").Value = 30
$wsPid.Range("P4").Value = 0.01
$wsPid.Range("Q4").Value = 1000
$wsPid.Range("R4").Value = 30
$wsPid.Range("S4").Value = 0

$wsPid.Range("A5").Value = 4
$wsPid.Range("B5").Value = 2.8
$wsPid.Range("C5").Value = 0
$wsPid.Range("D5").Value = 0.05
$wsPid.Range("E5").Value = 2.8
$wsPid.Range("F5").Value = 0
$wsPid.Range("G5").Value = 0.05
$wsPid.Range("H5").Value = 9
$wsPid.Range("I5").Value = 0.1
$wsPid.Range("J5").Value = 0
$wsPid.Range("K5").Value = 2800
$wsPid.Range("L5").Value = 35
$wsPid.Range("M5").Value = 0.01
$wsPid.Range("N5").Value = 2800
$wsPid.Range("O5").Value = 35
$wsPid.Range("P5").Value = 0.01
$wsPid.Range("Q5").Value = 1000
$wsPid.Range("R5").Value = 30
$wsPid.Range("S5").Value = 0

$wsPid.Range("A6").Value = 6

$wsPid.Range("A7").Value = 8
$wsPid.Range("B7").Value = 3.4
$wsPid.Range("C7").Value = 0
$wsPid.Range("D7").Value = 0.05
$wsPid.Range("E7").Value = 3.4
$wsPid.Range("F7").Value = 0
$wsPid.Range("G7").Value = 0.05
$wsPid.Range("H7").Value = 9
$wsPid.Range("I7").Value = 0.1
$wsPid.Range("J7").Value = 0
$wsPid.Range("K7").Value = 3200
$wsPid.Range("L7").Value = 50
$wsPid.Range("M7").Value = 0.01
$wsPid.Range("N7").Value = 3200
$wsPid.Range("O7").Value = 50
$wsPid.Range("P7").Value = 0.01
$wsPid.Range("Q7").Value = 2500
$wsPid.Range("R7").Value = 30
$wsPid.Range("S7").Value = 0

$wsPid.Range("A8").Value = 10

$wsPid.Range("A9").Value = 12
$wsPid.Range("B9").Value = 3.5
$wsPid.Range("C9").Value = 0
$wsPid.Range("D9").Value = 0.05
$wsPid.Range("E9").Value = 3.5
$wsPid.Range("F9").Value = 0
$wsPid.Range("G9").Value = 0.05
$wsPid.Range("H9").Value = 9
$wsPid.Range("I9").Value = 0.1
$wsPid.Range("J9").Value = 0
$wsPid.Range("K9").Value = 5500
$wsPid.Range("L9").Value = 100
$wsPid.Range("M9").Value = 0.01
$wsPid.Range("N9").Value = 5500
$wsPid.Range("O9").Value = 100
$wsPid.Range("P9").Value = 0.01
$wsPid.Range("Q9").Value = 2500
$wsPid.Range("R9").Value = 30

$wsPid.Range("A10").Value = 16

$wsPid.Range("A11").Value = 20

$wsPid.Range("A12").Value = 24

# --- view state: selections, scroll position, active sheet, zoom ----------
$ws1.Activate()
$ws1.Range("K63").Select()
$excel.ActiveWindow.Zoom = 95

$wsPid.Activate()
$wsPid.Range("T4").Select()
$excel.ActiveWindow.Zoom = 95

$wsExp.Activate()
$wsExp.Range("B5").Select()
$excel.ActiveWindow.Zoom = 95
